$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily row of data (row 84) following the existing pattern
$row = 84
$ws.Cells.Item($row, 1).Value = 46033
$ws.Cells.Item($row, 2).Value = 189
$ws.Cells.Item($row, 3).Value = 203
$ws.Cells.Item($row, 4).Value = 190

# Match the date number formatting used by the rest of column A
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
